$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Identifiers" column renamed to "Internal House Name" (same cell/position, O1)
$ws.Range("O1").Value = "Internal House Name"

# New trailing column "Tag /Band" header added in S1
$ws.Range("S1").Value = "Tag /Band"

# Widen the (renamed) column O to fit its new, longer header text
$ws.Columns("O").ColumnWidth = 14.52

# Move the selection/view over to the newly added column
$ws.Range("S1").Select() | Out-Null
